$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 5980
$ws.Range("I3").Value = 6233
$ws.Range("D4").Value = 1936
$ws.Range("I4").Value = 1433
$ws.Range("I5").Value = 580
$ws.Range("I6").Value = 7067
$ws.Range("D7").Value = 28126
$ws.Range("I7").Value = 21293

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 191
$ws.Range("I3").Value = 302
$ws.Range("I6").Value = 253
$ws.Range("I7").Value = 822

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 212

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 68
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 165
$ws.Range("I3").Value = 153
$ws.Range("I7").Value = 498

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("D2").Value = 235
$ws.Range("I2").Value = 165
$ws.Range("I6").Value = 154
$ws.Range("I7").Value = 666
$ws.Range("I8").Value = 1276
$ws.Range("I9").Value = 105
$ws.Range("I10").Value = 152
$ws.Range("I11").Value = 320
$ws.Range("I12").Value = 48
$ws.Range("I14").Value = 120
$ws.Range("I15").Value = 246
$ws.Range("I16").Value = 61
$ws.Range("I18").Value = 159
$ws.Range("I19").Value = 590
$ws.Range("I20").Value = 532
$ws.Range("I22").Value = 58
$ws.Range("I24").Value = 61
$ws.Range("I29").Value = 1317
$ws.Range("I31").Value = 212
$ws.Range("I33").Value = 963
$ws.Range("I40").Value = 39
$ws.Range("I47").Value = 148
$ws.Range("I48").Value = 286
$ws.Range("I49").Value = 146
$ws.Range("I50").Value = 105
$ws.Range("I52").Value = 463
$ws.Range("I55").Value = 234
$ws.Range("I57").Value = 85
$ws.Range("I63").Value = 76
$ws.Range("I65").Value = 498
$ws.Range("I67").Value = 822
$ws.Range("I73").Value = 196
$ws.Range("I76").Value = 304
$ws.Range("I78").Value = 290
$ws.Range("I79").Value = 603
$ws.Range("I83").Value = 459
$ws.Range("I84").Value = 184
$ws.Range("I85").Value = 971
$ws.Range("I88").Value = 194
$ws.Range("I90").Value = 261
$ws.Range("I95").Value = 325
$ws.Range("I97").Value = 181
$ws.Range("D101").Value = 28126
$ws.Range("I101").Value = 21293

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 166
$ws.Range("I7").Value = 459

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 112
$ws.Range("I7").Value = 325

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 214
$ws.Range("I3").Value = 366
$ws.Range("I6").Value = 303
$ws.Range("I7").Value = 963

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 454
$ws.Range("I6").Value = 367
$ws.Range("I7").Value = 1317

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 199
$ws.Range("I7").Value = 590

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I4").Value = 35
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 141
$ws.Range("I7").Value = 304

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 275
$ws.Range("I3").Value = 373
$ws.Range("I6").Value = 245
$ws.Range("I7").Value = 971

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I3").Value = 31
$ws.Range("I6").Value = 21

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I3").Value = 30
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 72
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 107
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 69
$ws.Range("I7").Value = 234

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I6").Value = 177
$ws.Range("I7").Value = 603

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 154
$ws.Range("I6").Value = 186
$ws.Range("I7").Value = 532

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 48
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 164
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 463

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I5").Value = 8
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 148

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 133
$ws.Range("I7").Value = 320

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 58
$ws.Range("D4").Value = 15
$ws.Range("I4").Value = 16
$ws.Range("I6").Value = 35
$ws.Range("D7").Value = 235
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 194

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I4").Value = 77
$ws.Range("I6").Value = 413
$ws.Range("I7").Value = 1276

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 88
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 261

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I2").Value = 23
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("I3").Value = 17
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 218
$ws.Range("I6").Value = 175
$ws.Range("I7").Value = 666

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 61
